$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 393
$ws.Range("I2").Value = 1031
$ws.Range("J2").Value = 4105
$ws.Range("K2").Value = 21
$ws.Range("L2").Value = 1112
$ws.Range("M2").Value = 75
$ws.Range("N2").Value = 706
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 20
$ws.Range("Q2").Value = 6
$ws.Range("R2").Value = 60
$ws.Range("S2").Value = 478
$ws.Range("T2").Value = 681
$ws.Range("V2").Value = 6427
$ws.Range("X2").Value = 6424
$ws.Range("Y2").Value = 8
$ws.Range("Z2").Value = 93
$ws.Range("AA2").Value = 43
